$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case " de "/" del "/" el "/" la "/" los " connectors in municipality/state names ---
$ws.Range("B4").Value = "Pabellón De Arteaga"
$ws.Range("B13").Value = "Amatenango De La Frontera"
$ws.Range("B31").Value = "Hidalgo Del Parral"
$ws.Range("B50").Value = "San Juan De Sabinas"
$ws.Range("B56").Value = "Villa De Álvarez"
$ws.Range("A58").Value = "Ciudad De México"
$ws.Range("B81").Value = "San Luis Del Cordero"
$ws.Range("A85").Value = "Estado De México"
$ws.Range("B88").Value = "Ecatepec De Morelos"
$ws.Range("B92").Value = "Naucalpan De Juárez"
$ws.Range("B97").Value = "Tlalnepantla De Baz"
$ws.Range("B112").Value = "Purísima Del Rincón"
$ws.Range("B114").Value = "Silao De La Victoria"
$ws.Range("B116").Value = "Valle De Santiago"
$ws.Range("B119").Value = "Ajuchitlán Del Progreso"
$ws.Range("B122").Value = "Atoyac De Álvarez"
$ws.Range("B123").Value = "Chilapa De Álvarez"
$ws.Range("B125").Value = "Coyuca De Catalán"
$ws.Range("B138").Value = "Pachuca De Soto"
$ws.Range("B141").Value = "Atotonilco El Alto"
$ws.Range("B142").Value = "Autlán De Navarro"
$ws.Range("B145").Value = "Encarnación De Díaz"
$ws.Range("B148").Value = "Lagos De Moreno"
$ws.Range("B151").Value = "San Diego De Alejandría"
$ws.Range("B153").Value = "San Juan De Los Lagos"
$ws.Range("B155").Value = "San Miguel El Alto"
$ws.Range("B157").Value = "Tamazula De Gordiano"
$ws.Range("B160").Value = "Tepatitlán De Morelos"
$ws.Range("B163").Value = "Unión De Tula"
$ws.Range("B166").Value = "Zapotlán El Grande"
$ws.Range("B200").Value = "Puente De Ixtla"
$ws.Range("B201").Value = "Tetela Del Volcán"
$ws.Range("B202").Value = "Tlaltizapán De Zapata"
$ws.Range("B207").Value = "Ixtlán Del Río"
$ws.Range("B209").Value = "Santa María Del Oro"
$ws.Range("B215").Value = "San Nicolás De Los Garza"
$ws.Range("B218").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B219").Value = "Chalcatongo De Hidalgo"
$ws.Range("B220").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B221").Value = "Ixtlán De Juárez"
$ws.Range("B222").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B223").Value = "Oaxaca De Juárez"
$ws.Range("B233").Value = "Zimatlán De Álvarez"
$ws.Range("B243").Value = "San Salvador El Verde"
$ws.Range("B252").Value = "Cadereyta De Montes"
$ws.Range("B255").Value = "Axtla De Terrazas"
$ws.Range("B258").Value = "Santa María Del Río"
$ws.Range("B285").Value = "Contla De Juan Cuamatzi"
$ws.Range("B289").Value = "Muñoz De Domingo Arenas"
$ws.Range("B293").Value = "Tepetitla De Lardizábal"
$ws.Range("B294").Value = "Tetla De La Solidaridad"
$ws.Range("B330").Value = "Nochistlán De Mejía"
$ws.Range("B331").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B333").Value = "Villa De Cos"

# --- Tiny floating point correction ---
$ws.Range("D287").Value = 0.09826589595375725

# --- Remove trailing footer/metadata rows beyond the data table (339:343, 476:480) ---
$ws.Range("A339:D343").ClearContents()
$ws.Range("A476:D480").ClearContents()
